# Auto-generated Excel COM-interop script to update FFXIV Leve profit data
# across the Pandaemonium Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Mirrors a scheduled market-data refresh: columns H-N (price / profit figures)
# are overwritten per-row; where the new NQ profit is undefined the cell is cleared.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4074.5833
$ws.Range("I64").Value = 3586
$ws.Range("J64").Value = 4888.8887
$ws.Range("K64").Value = 3586
$ws.Range("L64").Value = 4888.8887
$ws.Range("M64").Value = -3338
$ws.Range("N64").Value = -5384.8887

$ws.Range("H67").Value = 4074.5833
$ws.Range("I67").Value = 3586
$ws.Range("J67").Value = 4888.8887
$ws.Range("K67").Value = 3586
$ws.Range("L67").Value = 4888.8887
$ws.Range("M67").Value = -2728
$ws.Range("N67").Value = -6604.8887

$ws.Range("H98").Value = 1592.6923
$ws.Range("I98").Value = 1592.6923
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1592.6923
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -94.69229999999993

$ws.Range("H122").Value = 1592.6923
$ws.Range("I122").Value = 1592.6923
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4778.0769
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2328.0769
$ws.Range("N122").Value = -25900

$ws.Range("H137").Value = 3355.5
$ws.Range("I137").Value = 1684
$ws.Range("J137").Value = 6141.3335
$ws.Range("K137").Value = 5052
$ws.Range("L137").Value = 18424.0005
$ws.Range("M137").Value = -2502
$ws.Range("N137").Value = -23524.0005

$ws.Range("H138").Value = 3843.8206
$ws.Range("I138").Value = 2640.1875
$ws.Range("J138").Value = 4154.4355
$ws.Range("K138").Value = 7920.5625
$ws.Range("L138").Value = 12463.3065
$ws.Range("M138").Value = -2780.5625
$ws.Range("N138").Value = -22743.3065

$ws.Range("H139").Value = 74790
$ws.Range("J139").Value = 74790
$ws.Range("L139").Value = 74790
$ws.Range("N139").Value = -85070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6331.3267
$ws.Range("I32").Value = 4708.5
$ws.Range("J32").Value = 20612.2
$ws.Range("K32").Value = 4708.5
$ws.Range("L32").Value = 20612.2
$ws.Range("M32").Value = -4421.5
$ws.Range("N32").Value = -21186.2

$ws.Range("H92").Value = 38019.8
$ws.Range("J92").Value = 38019.8
$ws.Range("L92").Value = 38019.8
$ws.Range("N92").Value = -43011.8

$ws.Range("H129").Value = 34588.125
$ws.Range("J129").Value = 35570.855
$ws.Range("L129").Value = 35570.855
$ws.Range("N129").Value = -45570.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 56385.25
$ws.Range("J132").Value = 56385.25
$ws.Range("L132").Value = 56385.25
$ws.Range("N132").Value = -66505.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 112.6
$ws.Range("I7").Value = 128.28572
$ws.Range("J7").Value = 76
$ws.Range("K7").Value = 128.28572
$ws.Range("L7").Value = 76
$ws.Range("M7").Value = -15.28572
$ws.Range("N7").Value = -302

$ws.Range("H31").Value = 2678.814
$ws.Range("I31").Value = 1946.5278
$ws.Range("J31").Value = 6444.857
$ws.Range("K31").Value = 1946.5278
$ws.Range("L31").Value = 6444.857
$ws.Range("M31").Value = -1651.5278
$ws.Range("N31").Value = -7034.857

$ws.Range("H34").Value = 2678.814
$ws.Range("I34").Value = 1946.5278
$ws.Range("J34").Value = 6444.857
$ws.Range("K34").Value = 1946.5278
$ws.Range("L34").Value = 6444.857
$ws.Range("M34").Value = -1744.5278
$ws.Range("N34").Value = -6848.857

$ws.Range("H107").Value = 562.5
$ws.Range("I107").Value = 514.6
$ws.Range("J107").Value = 627.8182
$ws.Range("K107").Value = 514.6
$ws.Range("L107").Value = 627.8182
$ws.Range("M107").Value = 1405.4
$ws.Range("N107").Value = -4467.8182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 25.5
$ws.Range("I2").Value = 19.5
$ws.Range("J2").Value = 27.136364
$ws.Range("K2").Value = 117
$ws.Range("L2").Value = 162.818184
$ws.Range("M2").Value = -4
$ws.Range("N2").Value = -388.818184

$ws.Range("H34").Value = 3253.6667
$ws.Range("I34").Value = 573.25
$ws.Range("J34").Value = 5398
$ws.Range("K34").Value = 1719.75
$ws.Range("L34").Value = 16194
$ws.Range("M34").Value = -1635.75
$ws.Range("N34").Value = -16362

$ws.Range("H39").Value = 13535.714
$ws.Range("J39").Value = 14538.462
$ws.Range("L39").Value = 43615.386
$ws.Range("N39").Value = -44203.386

$ws.Range("H107").Value = 848.2041
$ws.Range("I107").Value = 299.6154
$ws.Range("J107").Value = 1046.3055
$ws.Range("K107").Value = 898.8462000000001
$ws.Range("L107").Value = 3138.9165
$ws.Range("M107").Value = 1021.1538
$ws.Range("N107").Value = -6978.916499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8305
$ws.Range("I122").Value = 12150
$ws.Range("J122").Value = 4460
$ws.Range("K122").Value = 36450
$ws.Range("L122").Value = 13380
$ws.Range("M122").Value = -34000
$ws.Range("N122").Value = -18280

$ws.Range("H139").Value = 57910
$ws.Range("J139").Value = 57910
$ws.Range("L139").Value = 57910
$ws.Range("N139").Value = -68190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4434
$ws.Range("I7").Value = 4550
$ws.Range("J7").Value = 4356.6665
$ws.Range("K7").Value = 4550
$ws.Range("L7").Value = 4356.6665
$ws.Range("M7").Value = -4438
$ws.Range("N7").Value = -4580.6665

$ws.Range("H122").Value = 5585.5
$ws.Range("I122").Value = 4331.1333
$ws.Range("J122").Value = 7937.4375
$ws.Range("K122").Value = 12993.3999
$ws.Range("L122").Value = 23812.3125
$ws.Range("M122").Value = -10543.3999
$ws.Range("N122").Value = -28712.3125

$ws.Range("H126").Value = 4434
$ws.Range("I126").Value = 4550
$ws.Range("J126").Value = 4356.6665
$ws.Range("K126").Value = 13650
$ws.Range("L126").Value = 13069.9995
$ws.Range("M126").Value = -11180
$ws.Range("N126").Value = -18009.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 35847.332
$ws.Range("J69").Value = 35847.332
$ws.Range("L69").Value = 35847.332
$ws.Range("N69").Value = -37345.332

$ws.Range("H72").Value = 35847.332
$ws.Range("J72").Value = 35847.332
$ws.Range("L72").Value = 107541.996
$ws.Range("N72").Value = -115029.996

$ws.Range("H105").Value = 70615
$ws.Range("J105").Value = 70615
$ws.Range("L105").Value = 70615
$ws.Range("N105").Value = -77603

$ws.Range("H113").Value = 1559.1818
$ws.Range("I113").Value = 993.1429000000001
$ws.Range("J113").Value = 2549.75
$ws.Range("K113").Value = 2979.4287
$ws.Range("L113").Value = 7649.25
$ws.Range("M113").Value = -809.4287000000004
$ws.Range("N113").Value = -11989.25

$ws.Range("H132").Value = 4609.8887
$ws.Range("I132").Value = 3915
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 11745
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -9215
$ws.Range("N132").Value = -23058.9995

$ws.Range("H138").Value = 40855.223
$ws.Range("J138").Value = 40855.223
$ws.Range("L138").Value = 40855.223
$ws.Range("N138").Value = -51135.223

